$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only the count (B) column changes
$ws.Range("B2").Value = 95693

# Rows 3 and 4 swap their species/location data (A, E, F, G, H, Q, R),
# while each gets its own new value in column B.
$ws.Range("A3").Value = 112182513
$ws.Range("B3").Value = 77388
$ws.Range("E3").Value = 6446
$ws.Range("F3").Value = "Kolflarnlav"
$ws.Range("G3").Value = "Carbonicola anthracophila"
$ws.Range("H3").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 366384
$ws.Range("R3").Value = 6863258

$ws.Range("A4").Value = 112182361
$ws.Range("B4").Value = 77636
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 366598
$ws.Range("R4").Value = 6863309

# Rows 5 and 6: only the count (B) column changes
$ws.Range("B5").Value = 77039
$ws.Range("B6").Value = 90812
